# Apply the evaluation-result updates to the QuantitativeMetrics sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 6: "Runtime without error" -> "no", with a note explaining why.
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "Missing initial redirect"

# Row 7: "Assertion validity" -> clear out the yes/note (no longer applicable).
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Row 12: refreshed Code BLEU score + breakdown string.
$ws.Range("B12").Value = 0.2767403357656455
$ws.Range("C12").Value = "{'codebleu': 0.2767403357656455, 'ngram_match_score': 0.13503301480439905, 'weighted_ngram_match_score': 0.1651517714816262, 'syntax_match_score': 0.532967032967033, 'dataflow_match_score': 0.27380952380952384}"

# Move the active selection from B6 to B7, matching the refreshed view state.
$ws.Range("B7").Select()
